# Auto-generated: apply cell value updates per diff
# Updates cryptocurrency price/volume table to latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.473.12"
$ws.Range("E2").Value = "  +4.04%  "
$ws.Range("D3").Value = "2.442.96"
$ws.Range("E3").Value = "  +3.37%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'556.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.36%  "
$ws.Range("D6").Value = "'138.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.572"
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = "  +3.64%  "
$ws.Range("E10").Value = "  +3.15%  "
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("D13").Value = "'24.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.67%  "
$ws.Range("D14").Value = "2.877.98"
$ws.Range("D15").Value = "60.368.02"
$ws.Range("E15").Value = "  +3.92%  "
$ws.Range("D16").Value = "'0.0000140"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("D17").Value = "2.449.21"
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("D18").Value = "'11.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.81%  "
$ws.Range("E19").Value = "  +3.11%  "
$ws.Range("D20").Value = "'335.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").Value = "'6.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'64.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.00%  "
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("D25").Value = "'8.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "0.0₃0792"
$ws.Range("E28").Value = "  +6.83%  "
$ws.Range("D29").Value = "'1.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.64%  "
$ws.Range("D30").Value = "'6.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").Value = "'170.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("E35").Value = "  +5.05%  "
$ws.Range("D36").Value = "'4.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "'1.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").Value = "'40.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'318.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.36%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.414"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.05%  "
$ws.Range("D42").Value = "'144.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("D44").Value = "'19.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.30%  "
$ws.Range("D45").Value = "'0.0964"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").Value = "'0.0525"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.04%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.574"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("B48").Value = "Polygon"
$ws.Range("C48").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D48").Value = "'0.404"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.69%  "
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D50").Value = "'18.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("D51").Value = "'11.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "
